$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.816.62"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.706.52"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").Value = "'0.9943"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").Value = "'317.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D7").Value = "'0.3920"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.4077"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "'1.498"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'54.44"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").Value = "'0.9935"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "'0.08817"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'26.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +11.34%  "
$ws.Range("D14").Value = "'7.490"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "'8.180"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'0.00001359"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").Value = "1.697.62"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'97.78"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'0.07163"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "'20.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.54%  "
$ws.Range("D21").Value = "'7.304"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("D22").Value = "'0.9950"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").Value = "24.797.49"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'3.031"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").Value = "'2.321"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "'23.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "'166.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "'6.031"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +17.03%  "
$ws.Range("D30").Value = "'8.538"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").Value = "'144.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.41%  "
$ws.Range("D32").Value = "1.883.17"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.184"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +11.51%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.08826"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "'1.079"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'7.280"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.06%  "
$ws.Range("D37").Value = "'0.03108"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2815"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.8531"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +11.42%  "
$ws.Range("D40").Value = "'10.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "'0.09213"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "'14.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "'1.479"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").Value = "'17.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.87%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.703"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.7479"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.46%  "
$ws.Range("D47").Value = "'4.263"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "'1.397"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.17%  "
$ws.Range("D49").Value = "'0.9944"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "'140.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "'0.08271"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.74%  "
